# Revert "Software review reconciliation changes"
# 1) Roll the cached datetimeFigureOut placeholder text back from 9/19/2017 to 11/3/2017
#    on the slide master and every slide layout.
# 2) Restore the older, longer USGS software-approval disclaimer text (and the
#    slightly taller/shifted textbox that holds it) on slide 2.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "11/3/2017"

# --- Every slide layout ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "11/3/2017"
}

# --- Slide 2: restore the disclaimer textbox text + geometry ---
$slide2 = $p.Slides.Item(2)
$disclaimer = $slide2.Shapes.Item(4)

$para1 = "This software has been approved for release by the U.S. Geological " + `
    "Survey (USGS). Although the software has been subjected to rigorous review, the " + `
    "USGS reserves the right to update the software as needed pursuant to further " + `
    "analysis and review. No warranty, expressed or implied, is made by the USGS or the " + `
    "U.S. Government as to the functionality of the software and related material nor shall " + `
    "the fact of release constitute any such warranty. Furthermore, the software is released " + `
    "on condition that neither the USGS nor the U.S. Government shall be held liable for any " + `
    "damages resulting from its authorized or"
$para2 = "unauthorized use."

$disclaimer.TextFrame.TextRange.Text = $para1 + [char]13 + $para2

$disclaimer.Left = 65.5220472440945
$disclaimer.Top = 385.0268503937008
$disclaimer.Width = 587.228188976378
$disclaimer.Height = 109.05472440944882
